$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("BP2").Value = 0.4227059363421588
$ws.Range("BQ2").Value = 0.06578473301828035
$ws.Range("BR2").Value = 0.1275153822730551
$ws.Range("BT2").Value = 0.04269178916126325
$ws.Range("BU2").Value = 0.06484576086461655
$ws.Range("BV2").Value = 0.1053691451084542
$ws.Range("BY2").Value = 0.06304419831828581
$ws.Range("BZ2").Value = 0.005957675852389435
$ws.Range("CB2").Value = 0.03890336911065064
$ws.Range("CC2").Value = 0.00795597562398411
$ws.Range("CD2").Value = 0.03669751311926647
$ws.Range("CI2").Value = 0.01012588600595163
$ws.Range("CR2").Value = 0.003794644870401417
$ws.Range("CT2").Value = 0.004607990331242381
$ws.Range("BP3").Value = 0.01640721081252504
$ws.Range("BQ3").Value = 0.3130748243644427
$ws.Range("BR3").Value = 0.2769117236781629
$ws.Range("BS3").Value = 0.08029580102114885
$ws.Range("BU3").Value = 0.1137928369315358
$ws.Range("BV3").Value = 0.01526593092326119
$ws.Range("BW3").Value = 0.05177677722200827
$ws.Range("BY3").Value = 0.004719985190078639
$ws.Range("BZ3").Value = 0.03340220612611006
$ws.Range("CC3").Value = 0.04700210166907482
$ws.Range("CD3").Value = 0.01021753106551277
$ws.Range("CE3").Value = 0.03713307099613899
$ws.Range("D4").Value = 0.04895182667398544
$ws.Range("E4").Value = 0.05151734615735443
$ws.Range("F4").Value = 0.03503350416128496
$ws.Range("H4").Value = 0.008705892759931071
$ws.Range("I4").Value = 0.003232438266502352
$ws.Range("K4").Value = 0.007694559960135596
$ws.Range("L4").Value = 0.05804362605112446
$ws.Range("M4").Value = 0.09396087758704909
$ws.Range("N4").Value = 0.03167887734882384
$ws.Range("P4").Value = 0.02517993074788459
$ws.Range("Q4").Value = 0.07136494053934149
$ws.Range("R4").Value = 0.09460081614613614
$ws.Range("S4").Value = 0.1196221764025613
$ws.Range("T4").Value = 0.09360085364101749
$ws.Range("U4").Value = 0.01505105183086493
$ws.Range("Y4").Value = 0.01703112177042816
$ws.Range("Z4").Value = 0.03039226896215791
$ws.Range("AB4").Value = 0.01045613837458238
$ws.Range("AC4").Value = 0.01144984248885231
$ws.Range("AD4").Value = 0.08398838992224139
$ws.Range("AE4").Value = 0.05552912717028268
$ws.Range("AF4").Value = 0.02079649589407237
$ws.Range("AG4").Value = 0.00175321620015875
$ws.Range("AI4").Value = 0.01036468094322692
$ws.Range("BP5").Value = 0.3720174631488403
$ws.Range("BQ5").Value = 0.09856167432464744
$ws.Range("BR5").Value = 0.1922017922463169
$ws.Range("BT5").Value = 0.0749929756845335
$ws.Range("BU5").Value = 0.06528056731492549
$ws.Range("BV5").Value = 0.03436211242205402
$ws.Range("BY5").Value = 0.03683042264633811
$ws.Range("CB5").Value = 0.04631611148075868
$ws.Range("CD5").Value = 0.07132536454217549
$ws.Range("CK5").Value = 0.008111516189410294
$ws.Range("BP6").Value = 0.003910122905935161
$ws.Range("BQ6").Value = 0.2546542562376047
$ws.Range("BR6").Value = 0.3026242938100749
$ws.Range("BS6").Value = 0.08164550971450804
$ws.Range("BU6").Value = 0.09043133490418438
$ws.Range("BW6").Value = 0.01483144740039066
$ws.Range("BY6").Value = 0.01363764945311798
$ws.Range("BZ6").Value = 0.02872881380570465
$ws.Range("CC6").Value = 0.07909232820890717
$ws.Range("CE6").Value = 0.1016483482845942
$ws.Range("CH6").Value = 0.002785089037388248
$ws.Range("CK6").Value = 0.01492916609074375
$ws.Range("CO6").Value = 0.0005975447083740931
$ws.Range("CP6").Value = 0.01048409543847212

$ws = $wb.Worksheets.Item(2)
$ws.Range("BP2").Value = 0.4227059363421588
$ws.Range("BQ2").Value = 0.4884906693604392
$ws.Range("BR2").Value = 0.6160060516334942
$ws.Range("BS2").Value = 0.6160060516334942
$ws.Range("BT2").Value = 0.6586978407947575
$ws.Range("BU2").Value = 0.723543601659374
$ws.Range("BV2").Value = 0.8289127467678282
$ws.Range("BW2").Value = 0.8289127467678282
$ws.Range("BX2").Value = 0.8289127467678282
$ws.Range("BY2").Value = 0.891956945086114
$ws.Range("BZ2").Value = 0.8979146209385035
$ws.Range("CA2").Value = 0.8979146209385035
$ws.Range("CB2").Value = 0.9368179900491541
$ws.Range("CC2").Value = 0.9447739656731382
$ws.Range("CD2").Value = 0.9814714787924047
$ws.Range("CE2").Value = 0.9814714787924047
$ws.Range("CF2").Value = 0.9814714787924047
$ws.Range("CG2").Value = 0.9814714787924047
$ws.Range("CH2").Value = 0.9814714787924047
$ws.Range("CI2").Value = 0.9915973647983564
$ws.Range("CJ2").Value = 0.9915973647983564
$ws.Range("CK2").Value = 0.9915973647983564
$ws.Range("CL2").Value = 0.9915973647983564
$ws.Range("CM2").Value = 0.9915973647983564
$ws.Range("CN2").Value = 0.9915973647983564
$ws.Range("CO2").Value = 0.9915973647983564
$ws.Range("CP2").Value = 0.9915973647983564
$ws.Range("CQ2").Value = 0.9915973647983564
$ws.Range("CR2").Value = 0.9953920096687577
$ws.Range("CS2").Value = 0.9953920096687577
$ws.Range("BP3").Value = 0.01640721081252504
$ws.Range("BQ3").Value = 0.3294820351769678
$ws.Range("BR3").Value = 0.6063937588551307
$ws.Range("BS3").Value = 0.6866895598762796
$ws.Range("BT3").Value = 0.6866895598762796
$ws.Range("BU3").Value = 0.8004823968078154
$ws.Range("BV3").Value = 0.8157483277310766
$ws.Range("BW3").Value = 0.8675251049530849
$ws.Range("BX3").Value = 0.8675251049530849
$ws.Range("BY3").Value = 0.8722450901431635
$ws.Range("BZ3").Value = 0.9056472962692736
$ws.Range("CA3").Value = 0.9056472962692736
$ws.Range("CB3").Value = 0.9056472962692736
$ws.Range("CC3").Value = 0.9526493979383485
$ws.Range("CD3").Value = 0.9628669290038613
$ws.Range("D4").Value = 0.04895182667398544
$ws.Range("E4").Value = 0.1004691728313399
$ws.Range("F4").Value = 0.1355026769926248
$ws.Range("G4").Value = 0.1355026769926248
$ws.Range("H4").Value = 0.1442085697525559
$ws.Range("I4").Value = 0.1474410080190583
$ws.Range("J4").Value = 0.1474410080190583
$ws.Range("K4").Value = 0.1551355679791939
$ws.Range("L4").Value = 0.2131791940303183
$ws.Range("M4").Value = 0.3071400716173674
$ws.Range("N4").Value = 0.3388189489661912
$ws.Range("O4").Value = 0.3388189489661912
$ws.Range("P4").Value = 0.3639988797140759
$ws.Range("Q4").Value = 0.4353638202534174
$ws.Range("R4").Value = 0.5299646363995535
$ws.Range("S4").Value = 0.6495868128021147
$ws.Range("T4").Value = 0.7431876664431322
$ws.Range("U4").Value = 0.7582387182739971
$ws.Range("V4").Value = 0.7582387182739971
$ws.Range("W4").Value = 0.7582387182739971
$ws.Range("X4").Value = 0.7582387182739971
$ws.Range("Y4").Value = 0.7752698400444253
$ws.Range("Z4").Value = 0.8056621090065833
$ws.Range("AA4").Value = 0.8056621090065833
$ws.Range("AB4").Value = 0.8161182473811657
$ws.Range("AC4").Value = 0.827568089870018
$ws.Range("AD4").Value = 0.9115564797922593
$ws.Range("AE4").Value = 0.967085606962542
$ws.Range("AF4").Value = 0.9878821028566144
$ws.Range("AG4").Value = 0.9896353190567732
$ws.Range("AH4").Value = 0.9896353190567732
$ws.Range("BP5").Value = 0.3720174631488403
$ws.Range("BQ5").Value = 0.4705791374734877
$ws.Range("BR5").Value = 0.6627809297198046
$ws.Range("BS5").Value = 0.6627809297198046
$ws.Range("BT5").Value = 0.7377739054043382
$ws.Range("BU5").Value = 0.8030544727192637
$ws.Range("BV5").Value = 0.8374165851413177
$ws.Range("BW5").Value = 0.8374165851413177
$ws.Range("BX5").Value = 0.8374165851413177
$ws.Range("BY5").Value = 0.8742470077876557
$ws.Range("BZ5").Value = 0.8742470077876557
$ws.Range("CA5").Value = 0.8742470077876557
$ws.Range("CB5").Value = 0.9205631192684144
$ws.Range("CC5").Value = 0.9205631192684144
$ws.Range("CD5").Value = 0.9918884838105899
$ws.Range("CE5").Value = 0.9918884838105899
$ws.Range("CF5").Value = 0.9918884838105899
$ws.Range("CG5").Value = 0.9918884838105899
$ws.Range("CH5").Value = 0.9918884838105899
$ws.Range("CI5").Value = 0.9918884838105899
$ws.Range("CJ5").Value = 0.9918884838105899
$ws.Range("BP6").Value = 0.003910122905935161
$ws.Range("BQ6").Value = 0.2585643791435399
$ws.Range("BR6").Value = 0.5611886729536147
$ws.Range("BS6").Value = 0.6428341826681228
$ws.Range("BT6").Value = 0.6428341826681228
$ws.Range("BU6").Value = 0.7332655175723072
$ws.Range("BV6").Value = 0.7332655175723072
$ws.Range("BW6").Value = 0.7480969649726978
$ws.Range("BX6").Value = 0.7480969649726978
$ws.Range("BY6").Value = 0.7617346144258158
$ws.Range("BZ6").Value = 0.7904634282315205
$ws.Range("CA6").Value = 0.7904634282315205
$ws.Range("CB6").Value = 0.7904634282315205
$ws.Range("CC6").Value = 0.8695557564404277
$ws.Range("CD6").Value = 0.8695557564404277
$ws.Range("CE6").Value = 0.9712041047250218
$ws.Range("CF6").Value = 0.9712041047250218
$ws.Range("CG6").Value = 0.9712041047250218
$ws.Range("CH6").Value = 0.9739891937624101
$ws.Range("CI6").Value = 0.9739891937624101
$ws.Range("CJ6").Value = 0.9739891937624101
$ws.Range("CK6").Value = 0.9889183598531539
$ws.Range("CL6").Value = 0.9889183598531539
$ws.Range("CM6").Value = 0.9889183598531539
$ws.Range("CN6").Value = 0.9889183598531539
$ws.Range("CO6").Value = 0.989515904561528

$ws = $wb.Worksheets.Item(3)
$ws.Range("D2").Value = 69
$ws.Range("F2").Value = 0.6160060516334942
$ws.Range("G2").Value = 4
$ws.Range("F3").Value = 0.6063937588551307
$ws.Range("F4").Value = 0.5299646363995535
$ws.Range("D5").Value = 69
$ws.Range("F5").Value = 0.6627809297198046
$ws.Range("G5").Value = 4
$ws.Range("F6").Value = 0.5611886729536147

$ws = $wb.Worksheets.Item(4)
$ws.Range("D2").Value = 72
$ws.Range("F2").Value = 0.723543601659374
$ws.Range("G2").Value = 7
$ws.Range("D3").Value = 72
$ws.Range("F3").Value = 0.8004823968078154
$ws.Range("G3").Value = 6
$ws.Range("F4").Value = 0.7431876664431322
$ws.Range("D5").Value = 71
$ws.Range("F5").Value = 0.7377739054043382
$ws.Range("G5").Value = 6
$ws.Range("D6").Value = 72
$ws.Range("F6").Value = 0.7332655175723072
$ws.Range("G6").Value = 6

$ws = $wb.Worksheets.Item(5)
$ws.Range("D2").Value = 73
$ws.Range("F2").Value = 0.8289127467678282
$ws.Range("G2").Value = 8
$ws.Range("F3").Value = 0.8004823968078154
$ws.Range("D4").Value = 25
$ws.Range("F4").Value = 0.8056621090065833
$ws.Range("G4").Value = 23
$ws.Range("D5").Value = 72
$ws.Range("F5").Value = 0.8030544727192637
$ws.Range("G5").Value = 7
$ws.Range("D6").Value = 80
$ws.Range("F6").Value = 0.8695557564404277
$ws.Range("G6").Value = 14

$ws = $wb.Worksheets.Item(6)
$ws.Range("D2").Value = 79
$ws.Range("F2").Value = 0.9368179900491541
$ws.Range("G2").Value = 14
$ws.Range("D3").Value = 77
$ws.Range("F3").Value = 0.9056472962692736
$ws.Range("G3").Value = 11
$ws.Range("F4").Value = 0.9115564797922593
$ws.Range("D5").Value = 79
$ws.Range("F5").Value = 0.9205631192684144
$ws.Range("G5").Value = 14
$ws.Range("F6").Value = 0.9712041047250218
